# #5: insurance, claim, debt, investment done
# Rework the "保險" (insurance) worksheet (sheet7) so every row carries the
# full common schema (company, name, owner, property_category, category,
# date, legislator_name, legislator_id, source_file, index) instead of the
# old (company, name, owner, premium-text) layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("保險")

# ---- header row ----------------------------------------------------------
$ws.Cells.Item(1, 2).Value = "company"
$ws.Cells.Item(1, 3).Value = "name"
$ws.Cells.Item(1, 4).Value = "owner"
$ws.Cells.Item(1, 5).Value = "property_category"
$ws.Cells.Item(1, 6).Value = "category"
$ws.Cells.Item(1, 7).Value = "date"
$ws.Cells.Item(1, 8).Value = "legislator_name"
$ws.Cells.Item(1, 9).Value = "legislator_id"
$ws.Cells.Item(1, 10).Value = "source_file"
$ws.Cells.Item(1, 11).Value = "index"

# ---- data rows ------------------------------------------------------------
# columns: A index, B company, C name, D owner, E property_category(=insurance)
#          F category(=normal), G date, H legislator_name, I legislator_id,
#          J source_file, K index (mirrors A)
$rows = @(
    @{ A = 109; B = "法國巴黎人壽"; C = "致富100富甲天下外幣變額保險"; D = "周桂香" },
    @{ A = 110; B = "南山人壽";     C = "增鑫動養老保險";             D = "周桂香" },
    @{ A = 111; B = "南山人壽";     C = "312還本終身保險勝";          D = "周桂香" },
    @{ A = 112; B = "南山人壽";     C = "美年發外幣增額還本終身保險"; D = "周桂香" },
    @{ A = 113; B = "南山人壽";     C = "312還本終身保險勝";          D = "周桂香" },
    @{ A = 114; B = "南山人壽";     C = "快樂兒童外幣增額還本終身保險"; D = "周桂香" },
    @{ A = 115; B = "南山人壽";     C = "美鑫外幣中身分紅壽險";       D = "周桂香" },
    @{ A = 116; B = "南山人壽";     C = "美鑫外幣中身分紅壽險";       D = "周桂香" },
    @{ A = 117; B = "南山人壽";     C = "伴我一生變額壽險";           D = "許智傑" },
    @{ A = 118; B = "南山人壽";     C = "美寶外幣終身分紅保險";       D = "周桂香" },
    @{ A = 119; B = "康健人壽";     C = "金準變額萬能壽險";           D = "周桂香" },
    @{ A = 120; B = "康健人壽";     C = "金準變額萬能壽險";           D = "周桂香" },
    @{ A = 121; B = "法國巴黎人壽"; C = "外幣變額年金保險";           D = "周桂香" },
    @{ A = 122; B = "富邦人壽";     C = "鑫添財萬能終身壽險";         D = "周桂香" },
    @{ A = 123; B = "國泰人壽";     C = "添美盛美元終身保險";         D = "周桂香" },
    @{ A = 124; B = "南山人壽";     C = "康祥一生終身保險";           D = "周桂香" },
    @{ A = 125; B = "南山人壽";     C = "鑫利年年外幣增額終身壽險";   D = "周桂香" },
    @{ A = 126; B = "南山人壽";     C = "康祥一生終身壽險";           D = "周桂香" },
    @{ A = 127; B = "南山人壽";     C = "康祥一生终身壽險";           D = "周桂香" },
    @{ A = 128; B = "南山人壽";     C = "新康祥終身壽險";             D = "周桂香" },
    @{ A = 129; B = "南山人壽";     C = "康祥一生終身壽險";           D = "周桂香" },
    @{ A = 130; B = "南山人壽";     C = "康祥一生終身壽險";           D = "周桂香" },
    @{ A = 132; B = "南山人壽";     C = "新20年期繳費增值分紅終身壽險"; D = "許智傑" },
    @{ A = 133; B = "南山人壽";     C = "新20年期缴費增值分紅終身壽險"; D = "許智傑" }
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = "insurance"
    $ws.Cells.Item($r, 6).Value = "normal"
    $ws.Cells.Item($r, 7).Value = "2013-12-11"
    $ws.Cells.Item($r, 8).Value = "許智傑"
    $ws.Cells.Item($r, 9).Value = 1750
    $ws.Cells.Item($r, 10).Value = "tmpd3cb1"
    $ws.Cells.Item($r, 11).Value = $row.A
    $r = $r + 1
}
